$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make sure column A cells are treated as plain text so the date-like
# strings ("05-10-2021", "06-10-2021") are not auto-converted to Excel
# date serial numbers, matching the existing "Serie" column formatting.
$ws.Range("A191:A192").NumberFormat = "@"

# Append two new daily rows of data (05-10-2021 and 06-10-2021)
$ws.Range("A191").Value = "05-10-2021"
$ws.Range("B191").Value = 1.69
$ws.Range("C191").Value = 2.16
$ws.Range("D191").Value = 3.12
$ws.Range("E191").Value = 2.4
$ws.Range("F191").Value = -0.62

$ws.Range("A192").Value = "06-10-2021"
$ws.Range("B192").Value = 1.91
$ws.Range("C192").Value = 2.29
$ws.Range("D192").Value = 3.12
$ws.Range("E192").Value = 2.4
$ws.Range("F192").Value = -0.59

# Restore default (General) number format on column A so the cells keep
# the same style as the rest of the "Serie" column (no explicit style).
$ws.Range("A191:A192").NumberFormat = "General"
